# Update "想去人数" (interested-people count) figures for four conge events
# that appear on both the "展览" sheet and the "全部类型" sheet.
#
#   15247 -> 15582
#   711   -> 713
#   643   -> 648
#   1625  -> 1630

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 15582
$wsExhibit.Range("F4").Value = 713
$wsExhibit.Range("F6").Value = 648
$wsExhibit.Range("F7").Value = 1630

# --- Sheet "全部类型" ---------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 15582
$wsAll.Range("F4").Value = 713
$wsAll.Range("F8").Value = 648
$wsAll.Range("F9").Value = 1630

$wb.Save()
